$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: Datos actualizados a 21 de Marzo de 2020 a las 14:16 -> Datos actualizados a 21 de Marzo de 2020 a las 14:46
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 21 de Marzo de 2020 a las 14:46"

# Row 7
$ws.Cells.Item(7, 2).Value = 21483
$ws.Cells.Item(7, 3).Value = 1635
$ws.Cells.Item(7, 5).Value = 21201
$ws.Cells.Item(7, 7).Value = 5
$ws.Cells.Item(7, 8).Value = 73

# Row 16
$ws.Cells.Item(16, 2).Value = 2785
$ws.Cells.Item(16, 3).Value = 136
$ws.Cells.Item(16, 5).Value = 2769

# Row 17
$ws.Cells.Item(17, 2).Value = 1999
$ws.Cells.Item(17, 3).Value = 40
$ws.Cells.Item(17, 5).Value = 1991
$ws.Cells.Item(17, 6).Value = 28

# Row 18
$ws.Cells.Item(18, 2).Value = 1764
$ws.Cells.Item(18, 3).Value = 125
$ws.Cells.Item(18, 5).Value = 1728

# Row 23
$ws.Cells.Item(23, 2).Value = 1072
$ws.Cells.Item(23, 3).Value = 144
$ws.Cells.Item(23, 5).Value = 1019

# Row 25
$ws.Cells.Item(25, 5).Value = 972
$ws.Cells.Item(25, 7).Value = 1
$ws.Cells.Item(25, 8).Value = 12

# Row 44
$ws.Cells.Item(44, 6).Value = 12

# Row 55
$ws.Cells.Item(55, 2).Value = 214
$ws.Cells.Item(55, 3).Value = 6
$ws.Cells.Item(55, 4).Value = 51
$ws.Cells.Item(55, 5).Value = 146

# Row 65: San Marino -> Emiratos Arabes Unidos
$ws.Cells.Item(65, 1).Value = "Emiratos Arabes Unidos"
$ws.Cells.Item(65, 2).Value = 153
$ws.Cells.Item(65, 3).Value = 13
$ws.Cells.Item(65, 4).Value = 31
$ws.Cells.Item(65, 5).Value = 120
$ws.Cells.Item(65, 6).Value = 2
$ws.Cells.Item(65, 8).Value = 2

# Row 66: Serbia -> San Marino
$ws.Cells.Item(66, 1).Value = "San Marino"
$ws.Cells.Item(66, 2).Value = 151
$ws.Cells.Item(66, 3).Value = 0
$ws.Cells.Item(66, 4).Value = 4
$ws.Cells.Item(66, 5).Value = 127
$ws.Cells.Item(66, 6).Value = 12
$ws.Cells.Item(66, 7).Value = 6
$ws.Cells.Item(66, 8).Value = 20

# Row 67: Bulgaria -> Serbia
$ws.Cells.Item(67, 1).Value = "Serbia"
$ws.Cells.Item(67, 2).Value = 149
$ws.Cells.Item(67, 3).Value = 14
$ws.Cells.Item(67, 4).Value = 2
$ws.Cells.Item(67, 5).Value = 146
$ws.Cells.Item(67, 6).Value = 4
$ws.Cells.Item(67, 8).Value = 1

# Row 68: Emiratos Arabes Unidos -> Bulgaria
$ws.Cells.Item(68, 1).Value = "Bulgaria"
$ws.Cells.Item(68, 2).Value = 142
$ws.Cells.Item(68, 3).Value = 15
$ws.Cells.Item(68, 4).Value = 3
$ws.Cells.Item(68, 5).Value = 136
$ws.Cells.Item(68, 6).Value = 3
$ws.Cells.Item(68, 8).Value = 3

# Row 74
$ws.Cells.Item(74, 5).Value = 51
$ws.Cells.Item(74, 7).Value = 1
$ws.Cells.Item(74, 8).Value = 12

# Row 81
$ws.Cells.Item(81, 2).Value = 85
$ws.Cells.Item(81, 5).Value = 84

# Row 83
$ws.Cells.Item(83, 2).Value = 77
$ws.Cells.Item(83, 3).Value = 4
$ws.Cells.Item(83, 5).Value = 74

# Row 113: Cuba -> Nigeria
$ws.Cells.Item(113, 1).Value = "Nigeria"
$ws.Cells.Item(113, 2).Value = 22
$ws.Cells.Item(113, 3).Value = 10
$ws.Cells.Item(113, 4).Value = 1
$ws.Cells.Item(113, 5).Value = 21
$ws.Cells.Item(113, 8).Value = 0

# Row 114: Ghana -> Cuba
$ws.Cells.Item(114, 1).Value = "Cuba"
$ws.Cells.Item(114, 2).Value = 21
$ws.Cells.Item(114, 3).Value = 0
$ws.Cells.Item(114, 5).Value = 20
$ws.Cells.Item(114, 8).Value = 1

# Row 115: Bolivia -> Ghana
$ws.Cells.Item(115, 1).Value = "Ghana"

# Row 116: Jamaica -> Bolivia
$ws.Cells.Item(116, 1).Value = "Bolivia"
$ws.Cells.Item(116, 3).Value = 3
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 19
$ws.Cells.Item(116, 8).Value = 0

# Row 117: Paraguay -> Jamaica
$ws.Cells.Item(117, 1).Value = "Jamaica"
$ws.Cells.Item(117, 2).Value = 19
$ws.Cells.Item(117, 4).Value = 2
$ws.Cells.Item(117, 5).Value = 16
$ws.Cells.Item(117, 6).Value = 0
$ws.Cells.Item(117, 7).Value = 0

# Row 118: Ruanda -> Paraguay
$ws.Cells.Item(118, 1).Value = "Paraguay"
$ws.Cells.Item(118, 2).Value = 18
$ws.Cells.Item(118, 6).Value = 1
$ws.Cells.Item(118, 7).Value = 1
$ws.Cells.Item(118, 8).Value = 1

# Row 119: Macao -> Ruanda
$ws.Cells.Item(119, 1).Value = "Ruanda"
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 17

# Row 120: Guayana Francesa -> Macao
$ws.Cells.Item(120, 1).Value = "Macao"
$ws.Cells.Item(120, 2).Value = 17
$ws.Cells.Item(120, 4).Value = 10
$ws.Cells.Item(120, 5).Value = 7

# Row 121: Polinesia Francesa -> Togo
$ws.Cells.Item(121, 1).Value = "Togo"
$ws.Cells.Item(121, 2).Value = 16
$ws.Cells.Item(121, 3).Value = 7
$ws.Cells.Item(121, 5).Value = 16

# Row 123: Kirguistan -> Guayana Francesa
$ws.Cells.Item(123, 1).Value = "Guayana Francesa"
$ws.Cells.Item(123, 2).Value = 15
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 5).Value = 15

# Row 124: Puerto Rico -> Polinesia Francesa
$ws.Cells.Item(124, 1).Value = "Polinesia Francesa"
$ws.Cells.Item(124, 2).Value = 15
$ws.Cells.Item(124, 3).Value = 4
$ws.Cells.Item(124, 5).Value = 15

# Row 125: Montenegro -> Puerto Rico
$ws.Cells.Item(125, 1).Value = "Puerto Rico"

# Row 126: Costa de Marfil -> Montenegro
$ws.Cells.Item(126, 1).Value = "Montenegro"
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 14

# Row 127: Mauricio -> Kirguistan
$ws.Cells.Item(127, 1).Value = "Kirguistan"
$ws.Cells.Item(127, 3).Value = 8
$ws.Cells.Item(127, 5).Value = 14
$ws.Cells.Item(127, 7).Value = 0
$ws.Cells.Item(127, 8).Value = 0

# Row 128: Guatemala -> Costa de Marfil
$ws.Cells.Item(128, 1).Value = "Costa de Marfil"
$ws.Cells.Item(128, 2).Value = 14
$ws.Cells.Item(128, 3).Value = 5
$ws.Cells.Item(128, 4).Value = 1
$ws.Cells.Item(128, 5).Value = 13
$ws.Cells.Item(128, 8).Value = 0

# Row 129: Maldivas -> Mauricio
$ws.Cells.Item(129, 1).Value = "Mauricio"
$ws.Cells.Item(129, 2).Value = 14
$ws.Cells.Item(129, 3).Value = 2
$ws.Cells.Item(129, 4).Value = 0
$ws.Cells.Item(129, 5).Value = 13
$ws.Cells.Item(129, 7).Value = 1
$ws.Cells.Item(129, 8).Value = 1

# Row 130: Nigeria -> Guatemala
$ws.Cells.Item(130, 1).Value = "Guatemala"
$ws.Cells.Item(130, 2).Value = 13
$ws.Cells.Item(130, 3).Value = 1
$ws.Cells.Item(130, 4).Value = 0
$ws.Cells.Item(130, 5).Value = 12
$ws.Cells.Item(130, 8).Value = 1

# Row 131: Monaco -> Maldivas
$ws.Cells.Item(131, 1).Value = "Maldivas"
$ws.Cells.Item(131, 2).Value = 13
$ws.Cells.Item(131, 4).Value = 2

# Row 132: Mongolia -> Monaco
$ws.Cells.Item(132, 1).Value = "Monaco"
$ws.Cells.Item(132, 2).Value = 11
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 5).Value = 11

# Row 133: Gibraltar -> Mongolia
$ws.Cells.Item(133, 1).Value = "Mongolia"
$ws.Cells.Item(133, 3).Value = 4
$ws.Cells.Item(133, 4).Value = 0
$ws.Cells.Item(133, 5).Value = 10

# Row 134: Etiopia -> Gibraltar
$ws.Cells.Item(134, 1).Value = "Gibraltar"
$ws.Cells.Item(134, 2).Value = 10
$ws.Cells.Item(134, 4).Value = 2
$ws.Cells.Item(134, 5).Value = 8

# Row 136: Togo -> Etiopia
$ws.Cells.Item(136, 1).Value = "Etiopia"

# Row 137: Seychelles -> Kenia
$ws.Cells.Item(137, 1).Value = "Kenia"

# Row 139: Kenia -> Seychelles
$ws.Cells.Item(139, 1).Value = "Seychelles"

# Row 140: Barbados -> Tanzania
$ws.Cells.Item(140, 1).Value = "Tanzania"

# Row 141: Tanzania -> Islas Virgenes de los Estados Unidos
$ws.Cells.Item(141, 1).Value = "Islas Virgenes de los Estados Unidos"
$ws.Cells.Item(141, 3).Value = 3

# Row 142: Guinea Ecuatorial -> Barbados
$ws.Cells.Item(142, 1).Value = "Barbados"

# Row 143: Islas Virgenes de los Estados Unidos -> Guinea Ecuatorial
$ws.Cells.Item(143, 1).Value = "Guinea Ecuatorial"
$ws.Cells.Item(143, 3).Value = 0

# Row 148: Bahamas -> San Martin (Parte Francesa)
$ws.Cells.Item(148, 1).Value = "San Martin (Parte Francesa)"

# Row 149: San Martin (Parte Francesa) -> Bahamas
$ws.Cells.Item(149, 1).Value = "Bahamas"

# Row 151: El Salvador -> San Bartolome
$ws.Cells.Item(151, 1).Value = "San Bartolome"
$ws.Cells.Item(151, 3).Value = 0

# Row 152: Namibia -> Zimbabue
$ws.Cells.Item(152, 1).Value = "Zimbabue"
$ws.Cells.Item(152, 3).Value = 2

# Row 153: San Bartolome -> El Salvador
$ws.Cells.Item(153, 1).Value = "El Salvador"
$ws.Cells.Item(153, 3).Value = 2

# Row 154: Liberia -> Madagascar
$ws.Cells.Item(154, 1).Value = "Madagascar"
$ws.Cells.Item(154, 3).Value = 0

# Row 155: Madagascar -> Congo
$ws.Cells.Item(155, 1).Value = "Congo"

# Row 156: Congo -> Republica de Africa Central
$ws.Cells.Item(156, 1).Value = "Republica de Africa Central"

# Row 157: Republica de Africa Central -> Namibia
$ws.Cells.Item(157, 1).Value = "Namibia"

# Row 158: Zimbabue -> Liberia
$ws.Cells.Item(158, 1).Value = "Liberia"
$ws.Cells.Item(158, 3).Value = 1

# Row 159: Islas Caimanes -> Cabo Verde
$ws.Cells.Item(159, 1).Value = "Cabo Verde"
$ws.Cells.Item(159, 3).Value = 2
$ws.Cells.Item(159, 5).Value = 3
$ws.Cells.Item(159, 8).Value = 0

# Row 161: Isla de Man -> Islas Caimanes
$ws.Cells.Item(161, 1).Value = "Islas Caimanes"
$ws.Cells.Item(161, 2).Value = 3
$ws.Cells.Item(161, 8).Value = 1

# Row 162: Groenlandia -> Angola
$ws.Cells.Item(162, 1).Value = "Angola"
$ws.Cells.Item(162, 3).Value = 1

# Row 164: Haiti -> Nicaragua
$ws.Cells.Item(164, 1).Value = "Nicaragua"

# Row 165: Zambia -> Benin
$ws.Cells.Item(165, 1).Value = "Benin"

# Row 166: Benin -> Guinea
$ws.Cells.Item(166, 1).Value = "Guinea"

# Row 167: Fiyi -> Mauritania
$ws.Cells.Item(167, 1).Value = "Mauritania"
$ws.Cells.Item(167, 3).Value = 0

# Row 168: Guinea -> Isla de Man
$ws.Cells.Item(168, 1).Value = "Isla de Man"

# Row 169: Angola -> Haiti
$ws.Cells.Item(169, 1).Value = "Haiti"
$ws.Cells.Item(169, 3).Value = 0

# Row 171: Santa Lucia -> Zambia
$ws.Cells.Item(171, 1).Value = "Zambia"

# Row 172: Nicaragua -> Groenlandia
$ws.Cells.Item(172, 1).Value = "Groenlandia"

# Row 173: Mauritania -> Fiyi
$ws.Cells.Item(173, 1).Value = "Fiyi"
$ws.Cells.Item(173, 3).Value = 1

# Row 174: Sudan -> Santa Lucia
$ws.Cells.Item(174, 1).Value = "Santa Lucia"
$ws.Cells.Item(174, 5).Value = 2
$ws.Cells.Item(174, 8).Value = 0

# Row 175: Cabo Verde -> Sudan
$ws.Cells.Item(175, 1).Value = "Sudan"
$ws.Cells.Item(175, 2).Value = 2
$ws.Cells.Item(175, 8).Value = 1

# Row 176: Republica de Yibuti -> San Vicente y las Granadinas
$ws.Cells.Item(176, 1).Value = "San Vicente y las Granadinas"

# Row 177: Santa Sede -> Gambia
$ws.Cells.Item(177, 1).Value = "Gambia"

# Row 178: San Martin (Parte Holandesa) -> Papua Nueva Guinea
$ws.Cells.Item(178, 1).Value = "Papua Nueva Guinea"

# Row 179: Niger -> Antigua y Barbuda
$ws.Cells.Item(179, 1).Value = "Antigua y Barbuda"

# Row 180: Timor Oriental -> Somalia
$ws.Cells.Item(180, 1).Value = "Somalia"
$ws.Cells.Item(180, 3).Value = 0

# Row 181: Somalia -> Suazilandia
$ws.Cells.Item(181, 1).Value = "Suazilandia"

# Row 182: Gambia -> Santa Sede
$ws.Cells.Item(182, 1).Value = "Santa Sede"

# Row 183: San Vicente y las Granadinas -> Montserrat
$ws.Cells.Item(183, 1).Value = "Montserrat"

# Row 184: Republica del Chad -> San Martin (Parte Holandesa)
$ws.Cells.Item(184, 1).Value = "San Martin (Parte Holandesa)"

# Row 185: Antigua y Barbuda -> Timor Oriental
$ws.Cells.Item(185, 1).Value = "Timor Oriental"
$ws.Cells.Item(185, 3).Value = 1

# Row 186: Papua Nueva Guinea -> Niger
$ws.Cells.Item(186, 1).Value = "Niger"

# Row 187: Suazilandia -> Republica de Yibuti
$ws.Cells.Item(187, 1).Value = "Republica de Yibuti"

# Row 188: Montserrat -> Republica del Chad
$ws.Cells.Item(188, 1).Value = "Republica del Chad"
